$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows in column I holding "ExcelReport-..." entries (Clinical, Economic,
# Quality of Life, Real-world Evidence) need their text changed so the
# "NewImportLogic_1 - Test_Automation_1" segment loses the spaces around
# the dash, becoming "NewImportLogic_1-Test_Automation_1".
$excelReportRows = @(3, 6, 9, 12)
foreach ($r in $excelReportRows) {
    $cell = $ws.Cells.Item($r, 9)
    $old = $cell.Value2
    $new = $old -replace "NewImportLogic_1 - Test_Automation_1", "NewImportLogic_1-Test_Automation_1"
    $cell.Value = $new
}

# Set explicit width for column I (new col definition in the diff).
# (76.5546875 "style units" ~= 75.72 characters of ColumnWidth.)
$ws.Columns.Item(9).ColumnWidth = 75.72

# Update the view: scroll so column H is the left-most visible column and
# select I13 only.
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("I13").Select()
